# ---------------------------------------------------------------------------
# codingHelper.xlsx edit: add a new "backtrack shortcut" worksheet summarising
# the backtracking-family LeetCode problems (subsets, permutations,
# combination sum, palindrome partitioning, ...), placed right after Sheet1
# and made the active sheet; also settle Sheet1's lingering selection.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# --- add the new sheet immediately after Sheet1 -----------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "backtrack shortcut"

# --- header row --------------------------------------------------------------
$newSheet.Range("A1").Value = 'Name'
$newSheet.Range("B1").Value = 'Difficulty'
$newSheet.Range("C1").Value = 'Solution'
$newSheet.Range("A2").Value = 'subsets'
$newSheet.Range("B2").Value = 'medium'
$newSheet.Range("C2").Value = 'Use backtrack. Empty arraylist. Add elements to it. Bactrack and then remove elements from it once all possibilities from that starting point are exhausted. Since empty set is a valid subset, make sure to add the empty arrayLst as well.'
$newSheet.Range("A3").Value = 'subsets without dups'
$newSheet.Range("B3").Value = 'medium'
$newSheet.Range("C3").Value = 'Same as subsets, for duplicate, for index greater than start idx, if char is same as previous, don’t do anything. Make sure to sort'
$newSheet.Range("A4").Value = 'permutations'
$newSheet.Range("B4").Value = 'medium'
$newSheet.Range("C4").Value = 'Use backtrack. Start with an empty arrayLsit and once its populated with all the elements(I.e len =  len of nums), then add it to final array. '
$newSheet.Range("A5").Value = 'permutations without dups'
$newSheet.Range("B5").Value = 'medium'
$newSheet.Range("C5").Value = 'Same as perm 1. Only diff is sort the array and maintian a used array. If this value has been used or the previous occurance of current value hasn''t been used yet, then don''t use this value, use it only after previous usage, otherwis, will repeat permutaitons.'
$newSheet.Range("A6").Value = 'combination sum'
$newSheet.Range("B6").Value = 'medium'
$newSheet.Range("C6").Value = 'Same as subsets, only difference is if the remainder value is 0, then add to final list, else don’t'
$newSheet.Range("A7").Value = 'combination sum without dups'
$newSheet.Range("B7").Value = 'medium'
$newSheet.Range("C7").Value = 'Same as subsets 2 with the condition used as combo sum'
$newSheet.Range("A8").Value = 'palindrome partitioning'
$newSheet.Range("B8").Value = 'medium'
$newSheet.Range("C8").Value = 'Use backtrack. Start with an empty arrayList, if there is a palindrome between start to I, then add it to the list. '
$newSheet.Range("A9").Value = 'Time complexity: '
$newSheet.Range("B9").Value = 'O(N* 2^N)'
$newSheet.Range("A10").Value = 'Space complexity:'
$newSheet.Range("B10").Value = 'O(N)'

# --- formatting ---------------------------------------------------------------
# Header row (A1:C1): bold + wrap, matching Sheet1's header style.
$sheet1.Range("A1").Copy()
$newSheet.Range("A1:C1").PasteSpecial(-4122)

# Rows 2-8 (the "medium" backtracking problems): themed fill + wrap text,
# matching the "Medium" rows on Sheet1.
$sheet1.Range("A2").Copy()
$newSheet.Range("A2:C8").PasteSpecial(-4122)

# Rows 9-10 (time/space complexity footnotes): plain wrap text, no fill.
$newSheet.Range("A9:B10").WrapText = $true

# --- column widths -------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 24
$newSheet.Columns.Item(2).ColumnWidth = 10
$newSheet.Columns.Item(3).ColumnWidth = 140.83

# --- selection / active sheet -------------------------------------------------
# Sheet1 no longer needs to keep the old A32 selection / scroll position.
$sheet1.Range("A1:C1").Select()

# Make the new sheet the active / selected tab, with C10 selected (matches
# where editing left off).
$newSheet.Range("C10").Select()
$newSheet.Select()

